# Adds the "Login GUI" section heading/spacer paragraphs after
# "Menu tap employees" and the full Reg010a-Reg010l GUI requirement
# list (plus trailing spacer paragraphs) after "Reg010:".
$d = $word.ActiveDocument

function Add-ParaAfter($para, [string]$xml) {
    # Insert a brand-new (empty) paragraph right after $para, then stamp its
    # final OOXML (runs / proofErr / lastRenderedPageBreak / ...) onto it via
    # InsertXML, which replaces the (empty) content of the range it is
    # called on.
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $newPara.Range.InsertXML($xml)
    return $newPara
}

function Find-ParaByText($doc, [string]$text) {
    $match = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            $match = $p
        }
    }
    return $match
}

# --- insert the 3 blank paragraphs + the new "Login GUI" heading
# --- paragraph right after "Menu tap employees" and before "Reg010:"
$anchor = Find-ParaByText $d "Menu tap employees"
if ($anchor -eq $null) { throw "Could not find anchor paragraph Menu tap employees" }
$anchor = Add-ParaAfter $anchor '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$anchor = Add-ParaAfter $anchor '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$anchor = Add-ParaAfter $anchor '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$anchor = Add-ParaAfter $anchor '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Login Graphical user interface (GUI) :</w:t></w:r></w:p>'

# --- locate "Reg010:" (the short heading paragraph, not Reg010a etc.)
$reg010 = Find-ParaByText $d "Reg010:"
if ($reg010 -eq $null) { throw "Could not find anchor paragraph Reg010:" }

# --- insert the Reg010a..Reg010l GUI-requirement paragraphs, the lone
# --- space paragraph, and the 2 trailing blank paragraphs after "Reg010:"
$a2 = $reg010
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reg010</w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t>welcome</w:t></w:r><w:r><w:t xml:space="preserve"> ma</w:t></w:r><w:r><w:t>ssage + username</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010b: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t xml:space="preserve">number of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>patince</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010c: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t xml:space="preserve">number of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>patince</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in the system</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010d: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t>number of doctors</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010e: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t>number of doctors in the system</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010f: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t>number of employees</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reg010</w:t></w:r><w:r><w:t>g</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>Label “</w:t></w:r><w:r><w:t>number of employees</w:t></w:r><w:r><w:t xml:space="preserve"> in the system</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reg010</w:t></w:r><w:r><w:t xml:space="preserve">h: </w:t></w:r><w:r><w:t>Button “</w:t></w:r><w:r><w:t>log out</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reg010i:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Button “</w:t></w:r><w:r><w:t>refresh stats</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Reg010j:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Menu tap </w:t></w:r><w:r><w:t>“</w:t></w:r><w:r><w:t>edit account</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010k: </w:t></w:r><w:r><w:t>Menu tap “</w:t></w:r><w:r><w:t>doctors</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Reg010l: </w:t></w:r><w:r><w:t>Menu tap “</w:t></w:r><w:r><w:t>employees</w:t></w:r><w:r><w:t>” visible: True</w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$a2 = Add-ParaAfter $a2 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'

Write-Output "paragraphs=$($d.Paragraphs.Count)"
